# Update the "Förändrad" (C column) date for every data row, and add a
# friendly-name second argument to the HYPERLINK() formulas on the first
# few rows that still have them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)  # column C ("Förändrad")
    $cCell.Value = 45186
}

# Columns that may contain HYPERLINK formulas: S(19) T(20) U(21) V(22) W(23) X(24) Y(25)
$hyperlinkCols = 19..25

for ($r = 2; $r -le $lastRow; $r++) {
    $nameCell = $ws.Cells.Item($r, 1)  # column A holds the "Beteckning" used as friendly name
    $name = $nameCell.Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ($formula -ne $null -and $formula -like '*HYPERLINK(*' -and $formula -notlike "*,*") {
            # formula looks like: =HYPERLINK("url")
            $newFormula = $formula -replace '\)\s*$', (', "' + $name + '")')
            $cell.Formula = $newFormula
        }
    }
}
